$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unmerge the A/B merged blocks first - columns A and B are merged across
# rows 2:18 and 19:24, and the boundary between the two groups is moving
# from row 18/19 to row 19/20. Writing to a non-top-left cell of a merged
# range is a no-op, so we must unmerge before editing cell values and
# re-merge with the new ranges afterwards.
$ws.Range("A2:A18").UnMerge()
$ws.Range("A19:A24").UnMerge()
$ws.Range("B2:B18").UnMerge()
$ws.Range("B19:B24").UnMerge()

# --- Row 12: now describes "Uploaded date" (was "File" / "Base64") ---
$ws.Range("E12").Value = "Uploaded date"
$ws.Range("F12").Value = ""
$ws.Range("G12").Value = "The date the document was uploaded to the application"
$ws.Range("H12").Value = "date"
$ws.Range("I12").Value = "MUST"

# --- Row 13: now describes "File" / "Base64" (was "Filename") ---
$ws.Range("F13").Value = "Base64"
$ws.Range("G13").Value = "Base64-encoded content of the file for inline file uploads"
$ws.Range("I13").Value = "MAY"

# --- Row 14: now describes "Filename" (was "MIME type") ---
$ws.Range("F14").Value = "Filename"
$ws.Range("G14").Value = "Name of the file being uploaded"
$ws.Range("I14").Value = "MUST"

# --- Row 15: now describes "MIME type" (was "File size") ---
$ws.Range("F15").Value = "MIME type"
$ws.Range("G15").Value = "The file's MIME type such as application/pdf or image/jpeg"
$ws.Range("H15").Value = "string"

# --- Row 16: now describes "Documents[] / File / File size" (was "Fee / Amount") ---
$ws.Range("D16").Value = "Documents[]"
$ws.Range("E16").Value = "File"
$ws.Range("F16").Value = "File size"
$ws.Range("G16").Value = "Size of the file in bytes that can be used to enforce limits"
$ws.Range("I16").Value = "MAY"

# --- Row 17: now describes "Fee / Amount" (was "Amount paid") ---
$ws.Range("E17").Value = "Amount"
$ws.Range("G17").Value = "The total amount due for the application fee"

# --- Row 18: now describes "Fee / Amount paid" (was "Transactions[]") ---
$ws.Range("E18").Value = "Amount paid"
$ws.Range("G18").Value = "The amount paid towards the application fee"
$ws.Range("H18").Value = "number"
$ws.Range("I18").Value = "MUST"

# --- Row 19: now describes "Fee / Transactions[]" (was "Advertisement location / Is advert in place") ---
$ws.Range("A19").Value = $null
$ws.Range("B19").Value = $null
$ws.Range("C19").Value = "Application"
$ws.Range("D19").Value = "Fee"
$ws.Range("E19").Value = "Transactions[]"
$ws.Range("G19").Value = "References to payments or financial transactions related to this application"
$ws.Range("H19").Value = "string"
$ws.Range("I19").Value = "MAY"

# --- Row 20: now starts "Advertisement location / Is advert in place" (was "Advert placed date") ---
$ws.Range("A20").Value = "Advertisement location"
$ws.Range("B20").Value = "Where the advertisement being applied to be built will be located"
$ws.Range("C20").Value = "Is advert in place"
$ws.Range("G20").Value = "Whether the advertisement is already in place"
$ws.Range("H20").Value = "boolean"
$ws.Range("I20").Value = "MUST"

# --- Row 21: now describes "Advert placed date" (was "Is replacement advert") ---
$ws.Range("C21").Value = "Advert placed date"
$ws.Range("G21").Value = "Date when the advertisement was placed (YYYY-MM-DD format)"
$ws.Range("H21").Value = "date"
$ws.Range("I21").Value = "MAY"

# --- Row 22: now describes "Is replacement advert" (was "Document reference[] / Reference") ---
$ws.Range("C22").Value = "Is replacement advert"
$ws.Range("D22").Value = ""
$ws.Range("G22").Value = "Whether this is a replacement advertisement"
$ws.Range("H22").Value = "boolean"

# --- Row 23: "Document reference[]" now describes "Reference" (was "Name") ---
$ws.Range("D23").Value = "Reference"
$ws.Range("G23").Value = "A unique reference for the data item"

# --- Row 79: datatype changed from "string" to "enum" ---
$ws.Range("H79").Value = "enum"

# --- Re-merge the A/B blocks using the new row boundary (19/20 instead of 18/19) ---
$ws.Range("A2:A19").Merge()
$ws.Range("A20:A24").Merge()
$ws.Range("B2:B19").Merge()
$ws.Range("B20:B24").Merge()
